$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title change: "Ghostly Spymaster" -> "Dark Dealings"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ghostly Spymaster", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dark Dealings", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Relocate the "_GoBack" bookmark from its old spot (near "Art tiles")
#    into the Overview paragraph, right after "...to sell t" and before
#    "o the black market...". Adding a bookmark named "_GoBack" replaces any
#    existing bookmark with that reserved name, so this both removes the old
#    one and creates the new one.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("to sell to the black market", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$splitPoint = $findRng.Start + [string]"to sell t".Length
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) Typo insertion: "guards" -> "gua``rds" (two literal backticks inserted)
#    within "Moving enemies like guards or dogs ...".
# ---------------------------------------------------------------------------
$replacement = 'gua``rds or dogs'
$d.Content.Find.Execute("guards or dogs", $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replacement, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Insert a new bullet after "Hub World/Main Menu- ..." paragraph:
#      - one blank ListParagraph (ind left=360)
#      - a new paragraph: "Level – Levels start in with you in a set
#        location, ..."
# ---------------------------------------------------------------------------
$hubPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Hub World/Main Menu*") {
        $hubPara = $p
        break
    }
}

if ($hubPara -ne $null) {
    $hubRange = $hubPara.Range
    $hubRange.InsertParagraphAfter()

    # find the freshly-inserted blank paragraph (right after the Hub World one)
    $blankPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -eq $hubRange.End) {
            $blankPara = $p
            break
        }
    }
    if ($blankPara -eq $null) {
        $blankPara = $hubPara.Next()
    }
    $blankRange = $blankPara.Range
    $blankRange.InsertParagraphAfter()

    $levelPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -eq $blankRange.End) {
            $levelPara = $p
            break
        }
    }
    if ($levelPara -eq $null) {
        $levelPara = $blankPara.Next()
    }
    $levelPara.Range.Text = "Level – Levels start in with you in a set location, the first phase you make your way towards the objective however you decide is best, once the objective is collected you then need to exit the level generally from your starting location. You can choose different paths to return if you like."
}

Write-Output "done"
